$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the command log (rows 136-155)
$newRows = @(
    ,@('187663574', '/start', 45124.94631944445)
    ,@('187663574', '/start', 45124.94635416667)
    ,@('187663574', 'Welcome курс | Адаптация', 45124.94658564815)
    ,@('187663574', 'Рассказывай!', 45124.94670138889)
    ,@('187663574', 'Оставить обращение', 45124.94673611111)
    ,@('187663574', 'Learning.telecom.kz | Техническая поддержка', 45124.94674768519)
    ,@('187663574', 'Обучение | Корпоративный Университет', 45124.9468287037)
    ,@('187663574', 'Служба поддержки “Нысана"', 45124.94685185186)
    ,@('187663574', 'Обратиться в службу комплаенс', 45124.946875)
    ,@('187663574', '/start', 45124.94689814815)
    ,@('187663574', 'База знаний', 45124.94692129629)
    ,@('187663574', 'База инструкций', 45124.94715277778)
    ,@('187663574', 'portal.telecom.kz | Инструкции', 45124.9471875)
    ,@('187663574', '/start', 45124.94725694445)
    ,@('187663574', 'Заполнить карточку БиОТ', 45124.94729166666)
    ,@('187663574', 'Опасный фактор/условие', 45124.94730324074)
    ,@('187663574', 'menu', 45124.94813657407)
    ,@('187663574', 'Часто задаваемые вопросы', 45124.9482175926)
    ,@('187663574', 'Демеу', 45124.94825231482)
    ,@('760906879', '/start', 45124.94930555556)
)

$startRow = 136
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $id = $newRows[$i][0]
    $cmd = $newRows[$i][1]
    $dt = $newRows[$i][2]
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $cmd
    $ws.Cells.Item($r, 3).Value = $dt
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

